$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-LatticeCell {
    param($cell, [string[]]$lines)
    $w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
    $runsXml = ""
    for ($i = 0; $i -lt $lines.Count; $i++) {
        $text = $lines[$i]
        $needsPreserve = ($text -ne $text.Trim())
        if ($needsPreserve) {
            $runsXml += "<w:t xml:space=`"preserve`">" + $text + "</w:t>"
        } else {
            $runsXml += "<w:t>" + $text + "</w:t>"
        }
        if ($i -lt ($lines.Count - 1)) {
            $runsXml += "<w:br/>"
        }
    }
    $xml = "<w:p xmlns:w=`"$w`"><w:r><w:rPr><w:sz w:val=`"32`"/></w:rPr>" + $runsXml + "</w:r></w:p>"
    $null = $cell.Range.InsertXML($xml)
}

Set-LatticeCell $t.Cell(1, 1) @("24 x 26", "  2    6", "  ----", "2|    |", "4|    |")
Set-LatticeCell $t.Cell(1, 2) @("16 x 72", "  7    2", "  ----", "1|    |", "6|    |")
Set-LatticeCell $t.Cell(1, 3) @("56 x 94", "  9    4", "  ----", "5|    |", "6|    |")
Set-LatticeCell $t.Cell(2, 1) @("71 x 55", "  5    5", "  ----", "7|    |", "1|    |")
Set-LatticeCell $t.Cell(2, 2) @("18 x 93", "  9    3", "  ----", "1|    |", "8|    |")
Set-LatticeCell $t.Cell(2, 3) @("29 x 97", "  9    7", "  ----", "2|    |", "9|    |")
Set-LatticeCell $t.Cell(3, 1) @("11 x 81", "  8    1", "  ----", "1|    |", "1|    |")
Set-LatticeCell $t.Cell(3, 2) @("44 x 83", "  8    3", "  ----", "4|    |", "4|    |")
Set-LatticeCell $t.Cell(3, 3) @("80 x 34", "  3    4", "  ----", "8|    |", "0|    |")
Set-LatticeCell $t.Cell(4, 1) @("27 x 89", "  8    9", "  ----", "2|    |", "7|    |")
Set-LatticeCell $t.Cell(4, 2) @("80 x 34", "  3    4", "  ----", "8|    |", "0|    |")
Set-LatticeCell $t.Cell(4, 3) @("74 x 79", "  7    9", "  ----", "7|    |", "4|    |")
Set-LatticeCell $t.Cell(5, 1) @("16 x 70", "  7    0", "  ----", "1|    |", "6|    |")
Set-LatticeCell $t.Cell(5, 2) @("31 x 56", "  5    6", "  ----", "3|    |", "1|    |")
Set-LatticeCell $t.Cell(5, 3) @("84 x 62", "  6    2", "  ----", "8|    |", "4|    |")
